# Auto-generated cell updates derived from the authoritative OOXML diff.
# Applies value changes (and the few cell adds/removes) per sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1375.7778
$ws.Range("J17").Value = 1375.7778
$ws.Range("L17").Value = 4127.3334
$ws.Range("N17").Value = -4463.3334
$ws.Range("H41").Value = 231
$ws.Range("J41").Value = 231
$ws.Range("L41").Value = 231
$ws.Range("N41").Value = -1111
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1381
$ws.Range("H48").Value = 1994
$ws.Range("J48").Value = 1994
$ws.Range("L48").Value = 5982
$ws.Range("N48").Value = -6566
$ws.Range("H52").Value = 2332.6667
$ws.Range("I52").Value = 998
$ws.Range("K52").Value = 2994
$ws.Range("M52").Value = -2834
$ws.Range("H56").Value = 1994
$ws.Range("J56").Value = 1994
$ws.Range("L56").Value = 5982
$ws.Range("N56").Value = -7050
$ws.Range("H60").Value = 500
$ws.Range("I60").Value = 500
$ws.Range("K60").Value = 1500
$ws.Range("M60").Value = -1016
$ws.Range("H62").Value = 5817.643
$ws.Range("I62").Value = 4883
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 4883
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -4259
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 5817.643
$ws.Range("I65").Value = 4883
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 24415
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -21295
$ws.Range("N65").Value = -43740
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H88").Value = 698.75
$ws.Range("I88").Value = 600
$ws.Range("J88").Value = 731.6667
$ws.Range("K88").Value = 600
$ws.Range("L88").Value = 731.6667
$ws.Range("M88").Value = -194
$ws.Range("N88").Value = -1543.6667
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H91").Value = 698.75
$ws.Range("I91").Value = 600
$ws.Range("J91").Value = 731.6667
$ws.Range("K91").Value = 600
$ws.Range("L91").Value = 731.6667
$ws.Range("M91").Value = 804
$ws.Range("N91").Value = -3539.6667
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 2375
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2375
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1834
$ws.Range("N100").Value = -4082
$ws.Range("H135").Value = 571.1667
$ws.Range("I135").Value = 571.1667
$ws.Range("K135").Value = 5140.5003
$ws.Range("M135").Value = -2605.5003
$ws.Range("H137").Value = 1980.7273
$ws.Range("I137").Value = 1258.5714
$ws.Range("J137").Value = 3244.5
$ws.Range("K137").Value = 3775.7142
$ws.Range("L137").Value = 9733.5
$ws.Range("M137").Value = -1225.7142
$ws.Range("N137").Value = -14833.5
$ws.Range("H138").Value = 1832.7142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2700
$ws.Range("I32").Value = 2405.2632
$ws.Range("K32").Value = 2405.2632
$ws.Range("M32").Value = -2118.2632
$ws.Range("H61").Value = 5204.4
$ws.Range("I61").Value = 4005.5
$ws.Range("K61").Value = 4005.5
$ws.Range("M61").Value = -3793.5
$ws.Range("H74").Value = 3864
$ws.Range("I74").Value = 3751.5
$ws.Range("K74").Value = 3751.5
$ws.Range("M74").Value = -2877.5
$ws.Range("H77").Value = 3864
$ws.Range("I77").Value = 3751.5
$ws.Range("K77").Value = 18757.5
$ws.Range("M77").Value = -14389.5
$ws.Range("H109").Value = 22222
$ws.Range("J109").Value = 22222
$ws.Range("L109").Value = 22222
$ws.Range("N109").Value = -24996
$ws.Range("H132").Value = 4503.75
$ws.Range("I132").Value = 4504.2856
$ws.Range("K132").Value = 13512.8568
$ws.Range("M132").Value = -10982.8568
$ws.Range("H136").Value = 5204.4
$ws.Range("I136").Value = 4005.5
$ws.Range("K136").Value = 12016.5
$ws.Range("M136").Value = -9466.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20400
$ws.Range("J92").Value = 20400
$ws.Range("L92").Value = 20400
$ws.Range("N92").Value = -25392
$ws.Range("H94").Value = 977.8570999999999
$ws.Range("J94").Value = 1065
$ws.Range("L94").Value = 1065
$ws.Range("N94").Value = -1967
$ws.Range("H102").Value = 25000
$ws.Range("I102").Value = 25000
$ws.Range("K102").Value = 25000
$ws.Range("M102").Value = -21755
$ws.Range("H134").Value = 5192.1333
$ws.Range("I134").Value = 5141.0713
$ws.Range("K134").Value = 15423.2139
$ws.Range("M134").Value = -12888.2139

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 10260.4
$ws.Range("I25").Value = 7092
$ws.Range("K25").Value = 7092
$ws.Range("M25").Value = -6918
$ws.Range("H31").Value = 2472.9092
$ws.Range("I31").Value = 2502.875
$ws.Range("K31").Value = 2502.875
$ws.Range("M31").Value = -2207.875
$ws.Range("H34").Value = 2472.9092
$ws.Range("I34").Value = 2502.875
$ws.Range("K34").Value = 2502.875
$ws.Range("M34").Value = -2300.875
$ws.Range("H50").Value = 20083.584
$ws.Range("I50").Value = 20111.445
$ws.Range("K50").Value = 20111.445
$ws.Range("M50").Value = -19486.445
$ws.Range("H56").Value = 14624.667
$ws.Range("I56").Value = 9437
$ws.Range("K56").Value = 9437
$ws.Range("M56").Value = -8592
$ws.Range("H58").Value = 2364
$ws.Range("I58").Value = 1460.4
$ws.Range("K58").Value = 1460.4
$ws.Range("M58").Value = -1257.4
$ws.Range("H59").Value = 28500.428
$ws.Range("I59").Value = 19752
$ws.Range("J59").Value = 31999.8
$ws.Range("K59").Value = 19752
$ws.Range("L59").Value = 31999.8
$ws.Range("M59").Value = -18607
$ws.Range("N59").Value = -34289.8
$ws.Range("H60").Value = 20847.928
$ws.Range("J60").Value = 19997.5
$ws.Range("L60").Value = 19997.5
$ws.Range("N60").Value = -21019.5
$ws.Range("H134").Value = 1274.6666
$ws.Range("I134").Value = 1274.6666
$ws.Range("K134").Value = 3823.9998
$ws.Range("M134").Value = -1288.9998
$ws.Range("H136").Value = 2364
$ws.Range("I136").Value = 1460.4
$ws.Range("K136").Value = 4381.200000000001
$ws.Range("M136").Value = -1831.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 75003
$ws.Range("J74").Value = 89999.664
$ws.Range("L74").Value = 269998.992
$ws.Range("N74").Value = -272120.992
$ws.Range("H77").Value = 75003
$ws.Range("J77").Value = 89999.664
$ws.Range("L77").Value = 809996.976
$ws.Range("N77").Value = -820604.976
$ws.Range("H122").Value = 1375.75
$ws.Range("J122").Value = 1002.5
$ws.Range("L122").Value = 9022.5
$ws.Range("N122").Value = -13922.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 170.53847
$ws.Range("I2").Value = 293
$ws.Range("J2").Value = 65.57143000000001
$ws.Range("K2").Value = 293
$ws.Range("L2").Value = 65.57143000000001
$ws.Range("M2").Value = -180
$ws.Range("N2").Value = -291.57143
$ws.Range("H5").Value = 49.5
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 49
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 63
$ws.Range("N5").Value = -274
$ws.Range("H111").Value = 70000
$ws.Range("J111").Value = 70000
$ws.Range("L111").Value = 70000
$ws.Range("N111").Value = -76134

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12412
$ws.Range("I7").Value = 11132.429
$ws.Range("K7").Value = 11132.429
$ws.Range("M7").Value = -11020.429
$ws.Range("H61").Value = 1936.5
$ws.Range("I61").Value = 2018.3334
$ws.Range("K61").Value = 2018.3334
$ws.Range("M61").Value = -1816.3334
$ws.Range("H98").Value = 15355
$ws.Range("J98").Value = 15355
$ws.Range("L98").Value = 15355
$ws.Range("N98").Value = -21345
$ws.Range("H113").Value = 1936.5
$ws.Range("I113").Value = 2018.3334
$ws.Range("K113").Value = 2018.3334
$ws.Range("M113").Value = 151.6666
$ws.Range("H126").Value = 12412
$ws.Range("I126").Value = 11132.429
$ws.Range("K126").Value = 33397.287
$ws.Range("M126").Value = -30927.287
$ws.Range("H132").Value = 4299
$ws.Range("I132").Value = 4995
$ws.Range("J132").Value = 3603
$ws.Range("K132").Value = 14985
$ws.Range("L132").Value = 10809
$ws.Range("M132").Value = -12455
$ws.Range("N132").Value = -15869

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 80000
$ws.Range("J121").Value = 80000
$ws.Range("L121").Value = 80000
$ws.Range("N121").Value = -83494
$ws.Range("H122").Value = 2599.6667
$ws.Range("I122").Value = 2599.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7799.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5349.000100000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 20435
$ws.Range("I126").Value = 14519.071
$ws.Range("K126").Value = 43557.213
$ws.Range("M126").Value = -41087.213

Write-Output "Applied 241 cell updates and 1 clear across 8 sheets."
